$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @('D2', '63.714.43'),
    @('E2', '  +3.08%  '),
    @('D3', '3.133.13'),
    @('E3', '  +1.42%  '),
    @('E4', '  -0.04%  '),
    @('D5', '591.11'),
    @('E5', '  +2.20%  '),
    @('D6', '146.10'),
    @('E6', '  +2.31%  '),
    @('E7', '  -0.02%  '),
    @('D8', '3.124.95'),
    @('E8', '  +1.54%  '),
    @('D9', '0.536'),
    @('E9', '  +1.62%  '),
    @('D10', '0.163'),
    @('E10', '  +17.22%  '),
    @('D11', '5.70'),
    @('E11', '  +3.94%  '),
    @('E12', '  +0.10%  '),
    @('D13', '0.0000255'),
    @('E13', '  +6.40%  '),
    @('D14', '36.18'),
    @('E14', '  +2.94%  '),
    @('E15', '  -0.51%  '),
    @('D16', '3.647.87'),
    @('E16', '  +1.53%  '),
    @('D17', '7.20'),
    @('D18', '63.627.86'),
    @('E18', '  +3.09%  '),
    @('D19', '3.127.32'),
    @('E19', '  +1.61%  '),
    @('D20', '466.05'),
    @('E20', '  +3.61%  '),
    @('D21', '14.27'),
    @('E21', '  +2.15%  '),
    @('D22', '0.737'),
    @('E22', '  +0.77%  '),
    @('E23', '  +1.12%  '),
    @('D24', '13.31'),
    @('E24', '  -2.19%  '),
    @('D25', '82.41'),
    @('E25', '  +0.26%  '),
    @('E26', '  -0.23%  '),
    @('D27', '8.76'),
    @('E27', '  +8.25%  '),
    @('D28', '2.72'),
    @('E28', '  +2.61%  '),
    @('D29', '2.23'),
    @('E30', '  -0.10%  '),
    @('E31', '  +1.50%  '),
    @('E32', '  +1.68%  '),
    @('E33', '  +0.09%  '),
    @('D34', '0.0₃0862'),
    @('E34', '  +8.11%  '),
    @('E35', '  +9.07%  '),
    @('E36', '  +1.90%  '),
    @('B37', 'dogwifhat'),
    @('C37', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'),
    @('D37', '3.37'),
    @('E37', '  +13.58%  '),
    @('B38', 'Filecoin'),
    @('C38', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'),
    @('D38', '6.14'),
    @('E38', '  +1.12%  '),
    @('D39', '50.90'),
    @('E39', '  +1.43%  '),
    @('D40', '449.30'),
    @('E40', '  +4.07%  '),
    @('D41', '8.76'),
    @('E41', '  -0.76%  '),
    @('E42', '  +0.58%  '),
    @('D43', '2.917.83'),
    @('E43', '  +4.61%  '),
    @('D44', '0.279'),
    @('E44', '  +3.76%  '),
    @('E45', '  +2.43%  '),
    @('D46', '2.17'),
    @('E46', '  +3.64%  '),
    @('D47', '125.00'),
    @('E47', '  -0.01%  '),
    @('D48', '35.04'),
    @('E48', '  -2.08%  '),
    @('E49', '  +0.01%  '),
    @('E50', '  +0.23%  '),
    @('D51', '24.72'),
    @('E51', '  +3.01%  ')
)

foreach ($edit in $edits) {
    $address = $edit[0]
    $text = $edit[1]
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}
